$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").NumberFormat = "General"
$ws.Range("A26").Value = 0
